# ------------------------------------------------------------------
# Commit: "[ADDITIONAL SCRAPING] added code to scrape more data about a
#          player's batting performance in a match, also updated the
#          excel sheets"
#
# Changes:
#   1. Add a new "Player Info" sheet before "ODI Batting" with player
#      bio data.
#   2. Add a new "ODI Batting Extra" sheet after "ODI Batting" with
#      additional per-match batting stats.
#   3. In "ODI Batting": rename column D header from MATCH_CARD_LINK to
#      MATCH_CODE, and replace every D-column URL value with just the
#      trailing MatchCode number (kept as text).
#   4. In "ODI Batting": clear the (empty) INNING_NUMBER cells for the
#      "did not bat" rows so they become true blank cells.
#
# NOTE: worksheet object variables in this COM shim resolve by
# *position*, not by identity -- inserting a sheet shifts what an
# already-held reference points to. So we always re-fetch sheets by
# name right before we use them, especially right after any
# Worksheets.Add() call.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------
# 1. New "Player Info" sheet, placed before "ODI Batting"
# --------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "'3779"
$playerInfo.Range("B2").Value = "Lokuge Dinesh Chandimal"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# Header styling, matching the bold/bordered/centered style used by the
# other header rows in this workbook (style index "1" in the source file).
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

# --------------------------------------------------------------
# 2. New "ODI Batting Extra" sheet, placed after "ODI Batting"
# --------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Add($null, $odiBatting)
$extra.Name = "ODI Batting Extra"

$extra = $wb.Worksheets.Item("ODI Batting Extra")
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Header styling, matching the bold/bordered/centered style used by the
# other header rows in this workbook (style index "1" in the source file).
$extra.Range("A1:F1").Font.Bold = $true
$extra.Range("A1:F1").HorizontalAlignment = -4108
$extra.Range("A1:F1").VerticalAlignment = -4160
$extra.Range("A1:F1").Borders.LineStyle = 1

$extraRows = @(
    @("4122", "4", "", "", "", "NO"),
    @("4124", "5", "0", "1", "20.36%", "NO"),
    @("4209", "", "", "", "", "NO"),
    @("4210", "3", "1", "0", "4.29%", "NO"),
    @("4211", "4", "2", "0", "22.67%", "NO"),
    @("4212", "3", "2", "0", "12.09%", "NO"),
    @("4215", "3", "6", "2", "21.86%", "NO"),
    @("4231", "5", "0", "0", "3.07%", "NO"),
    @("4232", "5", "0", "0", "1.01%", "NO"),
    @("4449", "5", "0", "0", "5.17%", "NO"),
    @("4450", "", "", "", "", "NO"),
    @("4451", "4", "2", "0", "5.84%", "NO"),
    @("4491", "2", "1", "0", "4.43%", "NO"),
    @("4521", "4", "4", "1", "25.00%", "YES"),
    @("4523", "", "", "", "", "NO"),
    @("4527", "4", "1", "0", "6.69%", "NO"),
    @("4603", "4", "1", "0", "3.75%", "NO"),
    @("4671", "", "", "", "", "NO"),
    @("4674", "3", "", "", "", "NO"),
    @("4675", "", "", "", "", "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = "'" + $row[0]
    if ($row[1] -ne "") {
        $extra.Cells.Item($r, 2).Value = [double]$row[1]
    }
    if ($row[2] -ne "") {
        $extra.Cells.Item($r, 3).Value = "'" + $row[2]
    }
    if ($row[3] -ne "") {
        $extra.Cells.Item($r, 4).Value = "'" + $row[3]
    }
    if ($row[4] -ne "") {
        $extra.Cells.Item($r, 5).Value = "'" + $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --------------------------------------------------------------
# 3. "ODI Batting": header + MATCH_CARD_LINK -> MATCH_CODE values
# --------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

for ($row = 2; $row -le 158; $row++) {
    $cell = $odiBatting.Cells.Item($row, 4)
    $url = $cell.Value2
    $idx = $url.IndexOf("MatchCode=")
    if ($idx -ge 0) {
        $code = $url.Substring($idx + 10)
        $cell.Value = "'" + $code
    }
}

# --------------------------------------------------------------
# 4. "ODI Batting": clear the empty INNING_NUMBER cells for the
#    "did not bat" rows so no placeholder cell remains.
# --------------------------------------------------------------
$emptyInningRows = @(39, 50, 53, 56, 59, 66, 67, 69, 94, 95, 107, 128, 139, 141, 157)
foreach ($row in $emptyInningRows) {
    $odiBatting.Cells.Item($row, 2).ClearContents()
}

# --------------------------------------------------------------
# Keep the first tab ("Player Info") as the active sheet, matching
# the original workbook's activeTab="0".
# --------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Activate()
